$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$content26 = @'
h1: What is my worth?
p: Tall men. Really tall lived on this planet long ago in history. They used to live for 1000 years. Some generations, really advanced, could find the location of angels looking at stars. Some built pyramids. Some build planes that fly. Some went to moon. A variety of men have lived and died. Some achieved great successes, secured huge lands and ruled the planet. Some died a failure. Among all these generations, who is it that really mattered to the Creator of this planet?
p.b-left: People who prostrated to Allah (swt) the most?
p.b-left: People who ruled with justice and eliminated injustices. Did they matter the most?
p.b-left: People who taught Quran and spread its word? Did they matter most? 
h3: Who was it that mattered?
p: <b>Ibrahim (as)</b>. 
p.note: I am going to try to picture the situation he faced, using today’s situation. Try to be with me.
p: The present day, we have a good normal life. We have secured financial stability. Our parents, wife and kids are living a fine life. Kids will get good education, parents will get good treatment in hospitals and wife will get good freedom. Everything is normal in any normal person’s life. 
p: Making it abnormal would be a chaos. It will be a down hill run and no coming back. <b>Raising voice for injustice happening somewhere out in the courts</b> would be abnormal. Trying to feed a kid for lifetime while<b> compromising my own kids studies</b> would be abnormal. <b>Selling all of our savings out and giving it in sadaqah</b> is abnormal. <b>Fighting to protect the underprivileged</b> in a normal life would be abnormal. 
p: Running down this path of abnormality is not recommended, not allowed. Dying in the path of Allah is not <b>recommended</b> by our relatives and friends. Shahadat (Martyrdom) is considered a loss to a young wife. Praying for it, or saying to achieve it out loud, is bad.
h3: Lets say, someone raises his voice for justice in an unjust society. What will happen?
p: He will be kicked out the 1st day from his home. He will run around for protection out in the open. No one will give him protection if he keeps acting “abnormal”.  
p: He will break down and fall, may be tomorrow or some days ahead. He will surrender to the norms of society. He will try to act a little normal that he is given a shelter and survive.
p: If he sticks to “<b>raising voice against injustice</b>” and “<b>does not settle for any less he determined to</b>”, he will be attacked harder by people in power. His own relatives first, followed by colleagues and neighbours will try to pull him down out of “<b>love for his wellbeing</b>”. 
h3: Still if he does not give up?
p: He will be taken to asylum cells. May be given electric shocks. He will be put behind bars. He will be made bedridden forcefully. He will be dented to fit in the definition of a<b> ‘Normal Life’</b>.
h3: Who was Ibrahim (as)?
p: Ibrahim (as) was somebody who after the electric shocks, the phase behind bars and the hard hitting; did not settle. He kept on running towards the fire, was thrown into it. He was made to question his own sanity. Just a dream, made him put his own son down for “straight to the point”, ”no questions asked” sacrifice. 
p: From being the <b>only 1</b> on this planet to believe in the might of Allah (swt), he made it to <b>1.8 Billion people</b> believing in the might of Creator. We face the stone, he moved. We face the city, he built. We love the sacrifice, he made. Our Prophet went into isolation to find Allah, Ibrahim found. He looked towards Allah (swt) to make the Kaaba, a Kibla, Ibrahim built. 
h3: How hard did Ibrahim (as) fight?
p: Some beaten, declared abnormal goes to Chief Justice of some state and questions his credibility. He will be humiliated, threatened and pushed further lower. Ibrahim (as) did not go to Chief Justice. He went straight to the ruler ‘Namrood’. Ibrahim (as) had neither power nor reputation to settle an argument with him. He straight on went into questioning his credibility. 
quote: Have you not considered the one who argued with Abraham about his Lord [merely] because Allah had given him kingship? When Abraham said, "My Lord is the one who gives life and causes death," he said, "I give life and cause death." Abraham said, "Indeed, Allah brings up the sun from the east, so bring it up from the west." So the disbeliever was overwhelmed [by astonishment], and Allah does not guide the wrongdoing people. <br>- Surah Baqarah verse 258
p: More sufferings his way. Did he give up? Did he not fall into fire for Allah (swt)? Did he not lay his own son a sacrifice? Did he not settle for nothing but truth?
p: Somebody did it for the love of Allah (swt). Where do I stand? After saying 5 prayers and writing articles and talking high. Still at 0. My worth? Like anybody’s worth, is 0 when it comes to comparing my superficial love with Ibrahim (as)’s love.
p: <b>May Allah (swt) accept our efforts and make us brave enough that we follow the path of Ibrahim (as). <span class=lavendar>May we make our Creator proud of us..</span> Amen </b>
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
'@
$ws.Range("D26").Value = $content26

$ayats27 = @'
Surah Baqarah, 265 - 275
'@
$tags27 = @'
Overcoming fear, Getting rich, Fixing financial issues, Al Sadiq, Al Amin, Marriage of Prophet (PBUH)
'@
$content27 = @'
h1: Getting really wealthy while following Quran
p.note: Guys, today I am going to speak, like I am speaking to my younger brother. I managed to become really rich in material by walking on this path. I really want we all live in an enlightened world where we are free of financial problems. 
p: We all agree that Muhammad (saw) was the most successful person ever on earth. He was having a good life. His wives were happy with him. His followers used to admire him. Till today he is quoted as the most influential person in history. His trade was successful. Battles were also a success. During Treaty of Hudaiybia, he was opposed by the closest to him, yet his sole decision-making lead entire state of Madina conquer Mecca. No blood was shed. People during invasion were feeling satisfied. Disbelievers from Mecca were also terrified at first, but later as Muslims won the city of Mecca, they were relieved. No one was harmed. People were guaranteed a safe life. Some were known for torturing muslims but they were forgiven. It was a delightful journey for Prophet Muhammad (saw) and all muslims.
p: Muhammad (saw) was really a successful person. He faced few hardships here and there. But he fought with courage and stood his ground. He lived a contented life and enjoyed the wealthiest life on earth. 
p: I honestly believe we should follow Muhammad (saw) to solve all our problems. Even if we want to get rich in material today, Muhammad (saw) was a fine tradesman to follow. His investments bore fruit and he won the wealthiest woman in Mecca, Hazrat Khadija (ra). 
h3: Marriage of Muhammad (saw) with Khadija (ra)
p: Muhammad (saw) was an illiterate man at the age of 25 when he married Hazrat Khadija (ra). Hazrat Khadija (ra) was a really wealthy woman at the time of her marriage. Her own caravan for trades to Syria surpassed all other caravans of Qurayesh. 
p: She was not used to travelling with her own caravans and preferred sending somebody else to represent her in her trade. In 595, Khadijah (ra) needed a co-worker for a transaction in Syria. She chose Muhammad ibn Abdullah for the trade in Syria. With the permission of Abu Talib ibn Muttalib, his uncle, he was sent to Syria with one of Khadijah's servants, Maysarah. 
p: The caravan after returning from Syria, got double the profit Khadija (ra) had expected. Her servant Maysarah told Khadija (ra) about the honesty, trust worthiness and truthfulness of Muhammad (saw). The success of caravan earned Muhammad (saw) the titles of “Al Sadiq” and “Al Amin”. Khadija (ra) admired the honesty of Muhammad (saw) and asked for his hand. They got married.
h3: Should we marry a wealthy lady to get rich?
p: Yes if we find one, we should. But what about people who are already married. Also what about people who do not have a choice but to marry someone from not-so-rich family. 
p: Muhammad (saw)’s life before Prophethood emphasises on two key aspects.
p.b-left: 1. He was truthful in his affairs.
p.b-left: 2. He was trust worthy that he would never break his promises.
p: Financial well-being is directly linked with our personal preferences. If we are not so trust worthy people, we might not get rich sooner. It will take a lot of time to convince people to invest in our ideas. However, if we become trustworthy and truthful, people will start listening to us. We will get good connections and we will get rich sooner. This is a very simple rule of life. <b>Honest shop keepers earn more then dishonest shopkeepers.</b>
h3: Honesty.. Is it that simple?
p: Honesty is a very difficult trait to adopt. Just like patience and discipline. All these good traits that make us rich are comparatively difficult to embrace. In order to become honest and trust worthy we need to work on some major aspects of our life.
p.b-left: 1. Make ourself purer 5 times a day. Start saying all 5 prayers.
p.b-left: 2. Write some baseline rules on a piece of paper and do not compromise on these rules.
p.b-left: 3. Control impulsive reaction to situations. Listen, digest and respond in a low voice maturely. 
p: These all things actually are one thing. It is to become a good person who lives a happy life. If we really want to get out of financial burden, we need to fix our own selves. Avoid following desires for some time, may be a month. Once we are on track, our desires will superimpose our habits. We will start living a purer and brighter day. Every morning will become better than the previous morning. 
p: Step by step we will become more richer in our life. We are not bad people ourselves. We are just afraid of loss while making bolder decisions. Less risks less fun. More risks more adventure.
h3: After achieving honest, invest boldly
p.note: This part is important, be with me..
p: Being honest will make us fearless. We develop a habit of talking truth and putting our reputation at stake more often. As we tame our fears, we become courageous to face our fears. Just like we tamed our fear of losing reputation by becoming honest and truthful, <span class=lavendar>we have to tame our next big fear; <b>”losing stability while investing money”.</b> </span>
p: The fear of losing stability in our life, stops us from investing in our bold ideas. These ideas are mostly million dollar ideas never acted upon due to fear of poverty. But as we start investing more, our investments start bearing fruits. Some times we lose a little, but the good relationships and trust we have developed through being honest, compliments our losses. It becomes easier to lose sometimes and win sometimes.
quote: Satan threatens you with poverty and orders you to immorality, while Allah promises you forgiveness from Him and bounty. And Allah is all-Encompassing and Knowing. <br> - Surah Baqarah verse 268
h3: How much time to get rich?
p: If we change today, in a month we will have better health. 2 months, better food on table. 3 months, better relationships. 4 months, better home. 5 months, better life. 1 or 2 years down this road, an ideal home where everything is perfect. Some say, this world never gets perfect. They are wrong. People who have better health, better relationships and a better car really do not live a bad life. 
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>.
'@
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 43856
$ws.Range("B27").NumberFormat = "d-mmm-yy"
$ws.Range("C27").Value = $ayats27
$ws.Range("F27").Value = $tags27
$ws.Range("D27").Value = $content27
$ws.Range("E27").Value = "Qasim Ali"
$ws.Rows.Item(27).RowHeight = 409.6

$ws.Range("D27").Select() | Out-Null
